$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap F:V content between row pairs (14,15), (42,43), (70,71)
function Swap-RowRange($ws, $r1, $r2, $colStart, $colEnd) {
    $range1 = $ws.Range("$colStart$r1`:$colEnd$r1")
    $range2 = $ws.Range("$colStart$r2`:$colEnd$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

Swap-RowRange $ws 14 15 "F" "V"
Swap-RowRange $ws 42 43 "F" "V"
Swap-RowRange $ws 70 71 "F" "V"
# Row 74
$ws.Range("A73:V73").Copy($ws.Range("A74:V74"))
$ws.Cells.Item(74,1).Value2 = 73
$ws.Cells.Item(74,2).Value2 = "portugal"
$ws.Cells.Item(74,3).Value2 = "liga-portugal"
$ws.Cells.Item(74,4).Value2 = "2023-2024"
$ws.Cells.Item(74,5).Value2 = 45226.88541666666
$ws.Cells.Item(74,6).Value2 = "Arouca"
$ws.Cells.Item(74,7).Value2 = 0
$ws.Cells.Item(74,8).Value2 = "Moreirense"
$ws.Cells.Item(74,9).Value2 = 1
$ws.Cells.Item(74,10).Value2 = 2.18
$ws.Cells.Item(74,11).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(74,12).Value2 = 2.68
$ws.Cells.Item(74,13).Value2 = "27/10/2023 21:06"
$ws.Cells.Item(74,14).Value2 = 3.43
$ws.Cells.Item(74,15).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(74,16).Value2 = 3.24
$ws.Cells.Item(74,17).Value2 = "27/10/2023 21:06"
$ws.Cells.Item(74,18).Value2 = 3.46
$ws.Cells.Item(74,19).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(74,20).Value2 = 2.9
$ws.Cells.Item(74,21).Value2 = "27/10/2023 21:06"
$ws.Cells.Item(74,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/arouca-moreirense/beWgZY7f/"

# Row 75
$ws.Range("A73:V73").Copy($ws.Range("A75:V75"))
$ws.Cells.Item(75,1).Value2 = 74
$ws.Cells.Item(75,2).Value2 = "portugal"
$ws.Cells.Item(75,3).Value2 = "liga-portugal"
$ws.Cells.Item(75,4).Value2 = "2023-2024"
$ws.Cells.Item(75,5).Value2 = 45227.6875
$ws.Cells.Item(75,6).Value2 = "Portimonense"
$ws.Cells.Item(75,7).Value2 = 1
$ws.Cells.Item(75,8).Value2 = "Estoril"
$ws.Cells.Item(75,9).Value2 = 0
$ws.Cells.Item(75,10).Value2 = 2.17
$ws.Cells.Item(75,11).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(75,12).Value2 = 2.52
$ws.Cells.Item(75,13).Value2 = "28/10/2023 16:22"
$ws.Cells.Item(75,14).Value2 = 3.4
$ws.Cells.Item(75,15).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(75,16).Value2 = 3.64
$ws.Cells.Item(75,17).Value2 = "28/10/2023 16:22"
$ws.Cells.Item(75,18).Value2 = 3.6
$ws.Cells.Item(75,19).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(75,20).Value2 = 2.81
$ws.Cells.Item(75,21).Value2 = "28/10/2023 16:23"
$ws.Cells.Item(75,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/portimonense-estoril/21GDx8En/"

# Row 76
$ws.Range("A73:V73").Copy($ws.Range("A76:V76"))
$ws.Cells.Item(76,1).Value2 = 75
$ws.Cells.Item(76,2).Value2 = "portugal"
$ws.Cells.Item(76,3).Value2 = "liga-portugal"
$ws.Cells.Item(76,4).Value2 = "2023-2024"
$ws.Cells.Item(76,5).Value2 = 45227.79166666666
$ws.Cells.Item(76,6).Value2 = "Benfica"
$ws.Cells.Item(76,7).Value2 = 1
$ws.Cells.Item(76,8).Value2 = "Casa Pia"
$ws.Cells.Item(76,9).Value2 = 1
$ws.Cells.Item(76,10).Value2 = 1.22
$ws.Cells.Item(76,11).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(76,12).Value2 = 1.22
$ws.Cells.Item(76,13).Value2 = "28/10/2023 18:55"
$ws.Cells.Item(76,14).Value2 = 7.32
$ws.Cells.Item(76,15).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(76,16).Value2 = 6.95
$ws.Cells.Item(76,17).Value2 = "28/10/2023 18:58"
$ws.Cells.Item(76,18).Value2 = 13.22
$ws.Cells.Item(76,19).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(76,20).Value2 = 14.17
$ws.Cells.Item(76,21).Value2 = "28/10/2023 18:58"
$ws.Cells.Item(76,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/benfica-casa-pia/GWtkzFhl/"

# Row 77
$ws.Range("A73:V73").Copy($ws.Range("A77:V77"))
$ws.Cells.Item(77,1).Value2 = 76
$ws.Cells.Item(77,2).Value2 = "portugal"
$ws.Cells.Item(77,3).Value2 = "liga-portugal"
$ws.Cells.Item(77,4).Value2 = "2023-2024"
$ws.Cells.Item(77,5).Value2 = 45227.79166666666
$ws.Cells.Item(77,6).Value2 = "Vitoria Guimaraes"
$ws.Cells.Item(77,7).Value2 = 5
$ws.Cells.Item(77,8).Value2 = "Chaves"
$ws.Cells.Item(77,9).Value2 = 0
$ws.Cells.Item(77,10).Value2 = 1.62
$ws.Cells.Item(77,11).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(77,12).Value2 = 1.69
$ws.Cells.Item(77,13).Value2 = "28/10/2023 18:58"
$ws.Cells.Item(77,14).Value2 = 4.26
$ws.Cells.Item(77,15).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(77,16).Value2 = 4
$ws.Cells.Item(77,17).Value2 = "28/10/2023 18:58"
$ws.Cells.Item(77,18).Value2 = 5.63
$ws.Cells.Item(77,19).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(77,20).Value2 = 5.29
$ws.Cells.Item(77,21).Value2 = "28/10/2023 18:58"
$ws.Cells.Item(77,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/vitoria-guimaraes-chaves/8vH9wlat/"

# Row 78
$ws.Range("A73:V73").Copy($ws.Range("A78:V78"))
$ws.Cells.Item(78,1).Value2 = 77
$ws.Cells.Item(78,2).Value2 = "portugal"
$ws.Cells.Item(78,3).Value2 = "liga-portugal"
$ws.Cells.Item(78,4).Value2 = "2023-2024"
$ws.Cells.Item(78,5).Value2 = 45227.89583333334
$ws.Cells.Item(78,6).Value2 = "Gil Vicente"
$ws.Cells.Item(78,7).Value2 = 3
$ws.Cells.Item(78,8).Value2 = "Braga"
$ws.Cells.Item(78,9).Value2 = 3
$ws.Cells.Item(78,10).Value2 = 4.14
$ws.Cells.Item(78,11).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(78,12).Value2 = 5.15
$ws.Cells.Item(78,13).Value2 = "28/10/2023 21:23"
$ws.Cells.Item(78,14).Value2 = 4.03
$ws.Cells.Item(78,15).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(78,16).Value2 = 4.81
$ws.Cells.Item(78,17).Value2 = "28/10/2023 21:23"
$ws.Cells.Item(78,18).Value2 = 1.85
$ws.Cells.Item(78,19).Value2 = "11/10/2023 14:42"
$ws.Cells.Item(78,20).Value2 = 1.59
$ws.Cells.Item(78,21).Value2 = "28/10/2023 21:23"
$ws.Cells.Item(78,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/gil-vicente-braga/WEgteGUD/"

# Row 79
$ws.Range("A73:V73").Copy($ws.Range("A79:V79"))
$ws.Cells.Item(79,1).Value2 = 78
$ws.Cells.Item(79,2).Value2 = "portugal"
$ws.Cells.Item(79,3).Value2 = "liga-portugal"
$ws.Cells.Item(79,4).Value2 = "2023-2024"
$ws.Cells.Item(79,5).Value2 = 45228.6875
$ws.Cells.Item(79,6).Value2 = "Rio Ave"
$ws.Cells.Item(79,7).Value2 = 3
$ws.Cells.Item(79,8).Value2 = "SC Farense"
$ws.Cells.Item(79,9).Value2 = 4
$ws.Cells.Item(79,10).Value2 = 2.13
$ws.Cells.Item(79,11).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(79,12).Value2 = 2.49
$ws.Cells.Item(79,13).Value2 = "29/10/2023 16:21"
$ws.Cells.Item(79,14).Value2 = 3.43
$ws.Cells.Item(79,15).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(79,16).Value2 = 3.38
$ws.Cells.Item(79,17).Value2 = "29/10/2023 16:21"
$ws.Cells.Item(79,18).Value2 = 3.69
$ws.Cells.Item(79,19).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(79,20).Value2 = 3.03
$ws.Cells.Item(79,21).Value2 = "29/10/2023 16:21"
$ws.Cells.Item(79,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/rio-ave-sc-farense/G4KHySTh/"

# Row 80
$ws.Range("A73:V73").Copy($ws.Range("A80:V80"))
$ws.Cells.Item(80,1).Value2 = 79
$ws.Cells.Item(80,2).Value2 = "portugal"
$ws.Cells.Item(80,3).Value2 = "liga-portugal"
$ws.Cells.Item(80,4).Value2 = "2023-2024"
$ws.Cells.Item(80,5).Value2 = 45228.79166666666
$ws.Cells.Item(80,6).Value2 = "Estrela"
$ws.Cells.Item(80,7).Value2 = 1
$ws.Cells.Item(80,8).Value2 = "Famalicao"
$ws.Cells.Item(80,9).Value2 = 0
$ws.Cells.Item(80,10).Value2 = 2.63
$ws.Cells.Item(80,11).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(80,12).Value2 = 2.61
$ws.Cells.Item(80,13).Value2 = "29/10/2023 18:56"
$ws.Cells.Item(80,14).Value2 = 3.12
$ws.Cells.Item(80,15).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(80,16).Value2 = 3.26
$ws.Cells.Item(80,17).Value2 = "29/10/2023 18:50"
$ws.Cells.Item(80,18).Value2 = 3.03
$ws.Cells.Item(80,19).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(80,20).Value2 = 2.97
$ws.Cells.Item(80,21).Value2 = "29/10/2023 18:56"
$ws.Cells.Item(80,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/estrela-da-amadora-famalicao/6LZ1XCx7/"

# Row 81
$ws.Range("A73:V73").Copy($ws.Range("A81:V81"))
$ws.Cells.Item(81,1).Value2 = 80
$ws.Cells.Item(81,2).Value2 = "portugal"
$ws.Cells.Item(81,3).Value2 = "liga-portugal"
$ws.Cells.Item(81,4).Value2 = "2023-2024"
$ws.Cells.Item(81,5).Value2 = 45228.89583333334
$ws.Cells.Item(81,6).Value2 = "Vizela"
$ws.Cells.Item(81,7).Value2 = 0
$ws.Cells.Item(81,8).Value2 = "FC Porto"
$ws.Cells.Item(81,9).Value2 = 2
$ws.Cells.Item(81,10).Value2 = 5.96
$ws.Cells.Item(81,11).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(81,12).Value2 = 8.02
$ws.Cells.Item(81,13).Value2 = "29/10/2023 21:29"
$ws.Cells.Item(81,14).Value2 = 4.53
$ws.Cells.Item(81,15).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(81,16).Value2 = 4.83
$ws.Cells.Item(81,17).Value2 = "29/10/2023 21:29"
$ws.Cells.Item(81,18).Value2 = 1.55
$ws.Cells.Item(81,19).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(81,20).Value2 = 1.43
$ws.Cells.Item(81,21).Value2 = "29/10/2023 21:22"
$ws.Cells.Item(81,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/vizela-fc-porto/pCY5WWiD/"

# Row 82
$ws.Range("A73:V73").Copy($ws.Range("A82:V82"))
$ws.Cells.Item(82,1).Value2 = 81
$ws.Cells.Item(82,2).Value2 = "portugal"
$ws.Cells.Item(82,3).Value2 = "liga-portugal"
$ws.Cells.Item(82,4).Value2 = "2023-2024"
$ws.Cells.Item(82,5).Value2 = 45229.88541666666
$ws.Cells.Item(82,6).Value2 = "Boavista"
$ws.Cells.Item(82,7).Value2 = 0
$ws.Cells.Item(82,8).Value2 = "Sporting CP"
$ws.Cells.Item(82,9).Value2 = 2
$ws.Cells.Item(82,10).Value2 = 5.32
$ws.Cells.Item(82,11).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(82,12).Value2 = 8.77
$ws.Cells.Item(82,13).Value2 = "30/10/2023 21:14"
$ws.Cells.Item(82,14).Value2 = 4.33
$ws.Cells.Item(82,15).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(82,16).Value2 = 5.15
$ws.Cells.Item(82,17).Value2 = "30/10/2023 21:14"
$ws.Cells.Item(82,18).Value2 = 1.63
$ws.Cells.Item(82,19).Value2 = "11/10/2023 14:43"
$ws.Cells.Item(82,20).Value2 = 1.38
$ws.Cells.Item(82,21).Value2 = "30/10/2023 21:06"
$ws.Cells.Item(82,22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal/boavista-sporting-lisbon/CbzbYhN0/"

